$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at H:I (shifts existing H:Q header data to J:S)
$ws.Range("H1:I1").EntireColumn.Insert() | Out-Null

# Fill the new header cells
$ws.Range("H1").Value = "Province"
$ws.Range("I1").Value = "TypeOfCar"

# Build the Tahoma/black font formatting once via a transient named style,
# apply it to both new header cells, then drop the named style again so we
# end up with just a plain extra cell format (no extra named cellStyle).
$st = $wb.Styles.Add("TempHeaderStyle")
$st.Font.Name = "Tahoma"
$st.Font.Color = 0
$ws.Range("H1:I1").Style = "TempHeaderStyle"
$wb.Styles.Item("TempHeaderStyle").Delete() | Out-Null

# Match the resulting selection state (whole column I selected, as after
# inserting/editing the new column)
$ws.Columns("I:I").Select() | Out-Null
